$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 100
$ws.Range("B4").Value = 8723.919156434198
$ws.Range("B5").Value = 89076.78952001187
$ws.Range("B7").Value = 1421.199252239083
$ws.Range("B9").Value = 2963.524179714889
$ws.Range("B10").Value = 158479.7600053143
$ws.Range("B11").Value = 0.07719403248579779
$ws.Range("B12").Value = 0.263145882107538
$ws.Range("B13").Value = 0.3500000000000028
$ws.Range("B14").Value = 0.9982733914346132
$ws.Range("B15").Value = 0.8962031735588429
